$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.10149999999999
$ws.Range("B21").Value = 5.810399999999992
$ws.Range("B23").Value = 5.656300000000001
$ws.Range("B25").Value = 5.875699999999994
